# Add a new "2020" data column (E) to the 5.4.1 indicator sheet, mirroring
# the existing "2015" column (D) formatting, plus one new number-format
# style for the "Parenting" 2020 figure, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy column D's cell formatting onto column E for the rows that keep
#     the same look (blank divider row, header year, and the totals row). ---
$ws.Range("D3:D6").Copy() | Out-Null
$ws.Range("E3:E6").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats

$ws.Range("D8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null      # -4122 = xlPasteFormats

$excel.CutCopyMode = $false

# --- New 2020 values ---
$ws.Range("E4").Value = 2020
$ws.Range("E5").Value = 11.5
$ws.Range("E6").Value = 2.6
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 0.3

# E7 ("Parenting" 2020 figure) gets its own style: same font/alignment as
# the rest of the data column, but with an explicit "0.0" number format.
$ws.Range("E7").NumberFormat = "0.0"

# --- Move the active selection, as recorded in the saved view state. ---
$ws.Range("B15").Select() | Out-Null
